$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the Normal style across the used range so a second cellXfs entry
# (fontId 0, applyFont=true) gets created and referenced by every cell.
$ws.Range("A1:C3").Style = "Normal"

# Fix the formulas so they correctly sum the contributing rate components
# (addresses the LEAC_plot_iter.py len(rate_table)-1 indexing bug).
$ws.Range("B2").Formula = "=0.1417+0.1502"
$ws.Range("C2").Formula = "=0.0644+0.1502"

# Move the selection on to the next row, matching the author's workflow.
$ws.Range("A4:D4").Select()
